# "Meer aanpassingen product backlog"
# Updates the "Product backlog" sheet: reword a couple of existing cells and
# add a new backlog row (row 7) describing the pump-control algorithm story.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product backlog")

# Row 2 — add the missing "Acceptatie criteria" text for the interface requirement
$ws.Range("D2").Value = "Een interface met bruikbare en relevante data die door iedereen gebruikt, bediend en bekeken kan worden"

# Row 4 — reworded task list and a small text tweak (trailing period added)
$ws.Range("C4").Value = "1. Geschikte pomp uitzoeken`n2. Testen en de beste plaatsing van de pompen kiezen`n3. Pompen implementeren in het schaalmodel"
$ws.Range("D4").Value = "De pompen brengen water in een uit de hier voor bestemde ruimtes."
$ws.Rows.Item(4).RowHeight = 43.5

# Row 7 — new backlog item about automatically steering the pumps from sensor data
$ws.Range("A7").Value = "6. Algoritme om de pompen aan te sturen op basis van sensordata."
$ws.Range("B7").Value = "Als systeemeigenaar wil ik dat het systeem uit zichzelf de boerderij stabiliseert op basis van de gemeten data, zodat deze niet door een mens geregeld hoeft te worden en altijd nauwkeurig zal zijn. "
$ws.Range("C7").Value = "1. Testen gebruik van pompen + sensor(en)`n2. Implementeren en tunen tot gewenst gebruik"
$ws.Range("D7").Value = "De pompen worden automatisch aangestuurd door het algoritme op basis van de sensordata"
$ws.Range("E7").Value = "M"

# A7:D7 are brand-new cells — match the wrap/left/top formatting used by the
# rest of the backlog table (same look as row 6's A:D cells).
$ws.Range("A7:D7").WrapText = $true
$ws.Range("A7:D7").HorizontalAlignment = -4131
$ws.Range("A7:D7").VerticalAlignment = -4160

$ws.Rows.Item(7).RowHeight = 58

# Move the view/selection to the newly added row
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C9").Select()

Write-Output "edits applied"
